# Generate Report for Handback
# Update timestamps / status produced by the handback report generation.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-01 20:19:00"
$wsOverview.Range("G4").Value = "2016-09-01 20:19:00"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-09-01 20:18:56"
$wsZhCn.Range("H4").Value = "2016-09-01 20:18:56"
$wsZhCn.Range("K2").Value = "2016-09-01 20:19:21"
$wsZhCn.Range("K4").Value = "2016-09-01 20:19:21"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-09-01 20:19:00"
$wsDeDe.Range("H4").Value = "2016-09-01 20:19:00"
$wsDeDe.Range("K2").Value = "2016-09-01 20:19:28"
$wsDeDe.Range("K4").Value = "2016-09-01 20:19:28"
